# Update the data rows 2-32 with the new "winter DST" values (data shifted by one week,
# and the day now has 31 half-hourly rows instead of 42 due to the DST change),
# then remove the now-unused trailing rows (33-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45954
$ws.Range("B2").Value = 3.36
$ws.Range("C2").Value = 9.511
$ws.Range("A3").Value = 45954.01041666666
$ws.Range("B3").Value = 1.637
$ws.Range("C3").Value = 0.984
$ws.Range("A4").Value = 45954.02083333334
$ws.Range("B4").Value = 0.079
$ws.Range("C4").Value = 1.922
$ws.Range("A5").Value = 45954.03125
$ws.Range("B5").Value = 2.171
$ws.Range("C5").Value = 0.567
$ws.Range("A6").Value = 45954.04166666666
$ws.Range("B6").Value = 5.691
$ws.Range("C6").Value = 0
$ws.Range("A7").Value = 45954.05208333334
$ws.Range("B7").Value = 13.826
$ws.Range("C7").Value = 0
$ws.Range("A8").Value = 45954.0625
$ws.Range("B8").Value = 15.817
$ws.Range("C8").Value = 0
$ws.Range("A9").Value = 45954.07291666666
$ws.Range("B9").Value = 16.297
$ws.Range("C9").Value = 0
$ws.Range("A10").Value = 45954.08333333334
$ws.Range("B10").Value = 34.057
$ws.Range("C10").Value = 0
$ws.Range("A11").Value = 45954.09375
$ws.Range("B11").Value = 38.432
$ws.Range("C11").Value = 0
$ws.Range("A12").Value = 45954.10416666666
$ws.Range("B12").Value = 19.766
$ws.Range("C12").Value = 0
$ws.Range("A13").Value = 45954.11458333334
$ws.Range("B13").Value = 28.14
$ws.Range("C13").Value = 0
$ws.Range("A14").Value = 45954.125
$ws.Range("B14").Value = 26.857
$ws.Range("C14").Value = 0
$ws.Range("A15").Value = 45954.13541666666
$ws.Range("B15").Value = 5.652
$ws.Range("C15").Value = 0.062
$ws.Range("A16").Value = 45954.14583333334
$ws.Range("B16").Value = 2.309
$ws.Range("C16").Value = 0.164
$ws.Range("A17").Value = 45954.15625
$ws.Range("B17").Value = 1.746
$ws.Range("C17").Value = 0.985
$ws.Range("A18").Value = 45954.16666666666
$ws.Range("B18").Value = 0.431
$ws.Range("C18").Value = 0.734
$ws.Range("A19").Value = 45954.17708333334
$ws.Range("B19").Value = 0.088
$ws.Range("C19").Value = 0.46
$ws.Range("A20").Value = 45954.1875
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 7.611
$ws.Range("A21").Value = 45954.19791666666
$ws.Range("B21").Value = 0.278
$ws.Range("C21").Value = 1.984
$ws.Range("A22").Value = 45954.20833333334
$ws.Range("B22").Value = 9.701
$ws.Range("C22").Value = 1.067
$ws.Range("A23").Value = 45954.21875
$ws.Range("B23").Value = 18.911
$ws.Range("C23").Value = 0
$ws.Range("A24").Value = 45954.22916666666
$ws.Range("B24").Value = 18.096
$ws.Range("C24").Value = 0
$ws.Range("A25").Value = 45954.23958333334
$ws.Range("B25").Value = 5.451
$ws.Range("C25").Value = 6.555
$ws.Range("A26").Value = 45954.25
$ws.Range("B26").Value = 1.682
$ws.Range("C26").Value = 20.564
$ws.Range("A27").Value = 45954.26041666666
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 35.693
$ws.Range("A28").Value = 45954.27083333334
$ws.Range("B28").Value = 0.007
$ws.Range("C28").Value = 7.558
$ws.Range("A29").Value = 45954.28125
$ws.Range("B29").Value = 0.145
$ws.Range("C29").Value = 3.717
$ws.Range("A30").Value = 45954.29166666666
$ws.Range("B30").Value = 2.281
$ws.Range("C30").Value = 1.157
$ws.Range("A31").Value = 45954.30208333334
$ws.Range("B31").Value = 2.412
$ws.Range("C31").Value = 0.863
$ws.Range("A32").Value = 45954.3125
$ws.Range("B32").Value = 0.02
$ws.Range("C32").Value = 1.745

# Remove the rows that are no longer part of the dataset (rows 33 through 43)
$ws.Range("A33:C43").EntireRow.Delete()
